# Apply updates to the "Resumo de Inscricoes" worksheet to reflect the
# latest enrollment counts (Inscritos / Pagos / Inscricoes homologadas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 4: Inscritos 14 -> 15
$ws.Range("E4").Value = 15

# Row 7: Inscritos 21 -> 22
$ws.Range("E7").Value = 22

# Row 15: Inscritos 70 -> 71
$ws.Range("E15").Value = 71

# Row 16: Inscritos 256 -> 259, Pagos 69 -> 70, Inscricoes homologadas 69 -> 70
$ws.Range("E16").Value = 259
$ws.Range("F16").Value = 70
$ws.Range("H16").Value = 70
